$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-10 Tuesday" "2024-12-11 Wednesday"

Replace-Text "699÷7=99, 6" "197÷5=39, 2"
Replace-Text "582÷3=194, 0" "649÷5=129, 4"
Replace-Text "527÷2=263, 1" "635÷2=317, 1"
Replace-Text "996÷5=199, 1" "727÷3=242, 1"
Replace-Text "121÷8=15, 1" "551÷9=61, 2"

Replace-Text "707÷3=235, 2" "854÷4=213, 2"
Replace-Text "611÷2=305, 1" "850÷7=121, 3"
Replace-Text "990÷9=110, 0" "557÷7=79, 4"
Replace-Text "488÷9=54, 2" "246÷6=41, 0"
Replace-Text "986÷8=123, 2" "586÷3=195, 1"

Replace-Text "308÷5=61, 3" "255÷6=42, 3"
Replace-Text "312÷3=104, 0" "909÷3=303, 0"
Replace-Text "838÷9=93, 1" "524÷4=131, 0"
Replace-Text "902÷5=180, 2" "604÷5=120, 4"
Replace-Text "145÷8=18, 1" "808÷9=89, 7"

Replace-Text "905÷5=181, 0" "790÷9=87, 7"
Replace-Text "628÷2=314, 0" "886÷8=110, 6"
Replace-Text "880÷7=125, 5" "628÷7=89, 5"
Replace-Text "589÷2=294, 1" "436÷2=218, 0"
Replace-Text "947÷2=473, 1" "449÷5=89, 4"

Replace-Text "428÷6=71, 2" "914÷5=182, 4"
Replace-Text "508÷6=84, 4" "991÷3=330, 1"
Replace-Text "529÷6=88, 1" "449÷4=112, 1"
Replace-Text "398÷5=79, 3" "724÷6=120, 4"
Replace-Text "295÷8=36, 7" "382÷8=47, 6"
